$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2030651340996169
$ws.Range("C2").Value = 0.5363984674329502
$ws.Range("J2").Value = 0.01149425287356322
$ws.Range("P2").Value = 0.1532567049808429
$ws.Range("S2").Value = 0.09578544061302682
$ws.Range("B3").Value = 0.01986754966887417
$ws.Range("C3").Value = 0.05298013245033113
$ws.Range("J3").Value = 0.01986754966887417
$ws.Range("P3").Value = 0.8013245033112583
$ws.Range("S3").Value = 0.1059602649006623
$ws.Range("J4").Value = 0.02941176470588235
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.04210526315789474
$ws.Range("D6").Value = 0.01578947368421053
$ws.Range("F6").Value = 0.07894736842105263
$ws.Range("J6").Value = 0.2368421052631579
$ws.Range("O6").Value = 0.01578947368421053
$ws.Range("Q6").Value = 0.1631578947368421
$ws.Range("R6").Value = 0.07368421052631578
$ws.Range("S6").Value = 0.3736842105263158
$ws.Range("B7").Value = 0.08947368421052632
$ws.Range("D7").Value = 0.01052631578947368
$ws.Range("F7").Value = 0.07368421052631578
$ws.Range("J7").Value = 0.04736842105263158
$ws.Range("O7").Value = 0.03157894736842105
$ws.Range("Q7").Value = 0.1473684210526316
$ws.Range("R7").Value = 0.1263157894736842
$ws.Range("S7").Value = 0.4736842105263158
$ws.Range("B8").Value = 0.09574468085106383
$ws.Range("D8").Value = 0.01914893617021277
$ws.Range("F8").Value = 0.05319148936170213
$ws.Range("J8").Value = 0.1063829787234043
$ws.Range("O8").Value = 0.03191489361702127
$ws.Range("Q8").Value = 0.1808510638297872
$ws.Range("R8").Value = 0.1148936170212766
$ws.Range("S8").Value = 0.3978723404255319
$ws.Range("B9").Value = 0.06956521739130435
$ws.Range("D9").Value = 0.02173913043478261
$ws.Range("F9").Value = 0.0391304347826087
$ws.Range("J9").Value = 0.0782608695652174
$ws.Range("O9").Value = 0.04347826086956522
$ws.Range("Q9").Value = 0.1739130434782609
$ws.Range("R9").Value = 0.1304347826086956
$ws.Range("S9").Value = 0.4434782608695652
$ws.Range("B10").Value = 0.1013100436681223
$ws.Range("D10").Value = 0.01397379912663755
$ws.Range("F10").Value = 0.05414847161572053
$ws.Range("J10").Value = 0.0777292576419214
$ws.Range("O10").Value = 0.01834061135371179
$ws.Range("Q10").Value = 0.211353711790393
$ws.Range("R10").Value = 0.1004366812227074
$ws.Range("S10").Value = 0.422707423580786
$ws.Range("G11").Value = 0.1283783783783784
$ws.Range("J11").Value = 0.08108108108108109
$ws.Range("K11").Value = 0.1722972972972973
$ws.Range("L11").Value = 0.581081081081081
$ws.Range("S11").Value = 0.03716216216216216
$ws.Range("G12").Value = 0.7443181818181818
$ws.Range("J12").Value = 0.1761363636363636
$ws.Range("K12").Value = 0.005681818181818182
$ws.Range("L12").Value = 0.02840909090909091
$ws.Range("S12").Value = 0.04545454545454546
$ws.Range("F13").Value = 0.02380952380952381
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2142857142857143
$ws.Range("S13").Value = 0.09523809523809523
$ws.Range("F15").Value = 0.01702127659574468
$ws.Range("H15").Value = 0.1659574468085106
$ws.Range("I15").Value = 0.07659574468085106
$ws.Range("J15").Value = 0.2851063829787234
$ws.Range("K15").Value = 0.06808510638297872
$ws.Range("M15").Value = 0.008510638297872341
$ws.Range("O15").Value = 0.04680851063829787
$ws.Range("S15").Value = 0.3319148936170213
$ws.Range("F16").Value = 0.02285714285714286
$ws.Range("H16").Value = 0.1828571428571429
$ws.Range("I16").Value = 0.06857142857142857
$ws.Range("J16").Value = 0.4342857142857143
$ws.Range("K16").Value = 0.08571428571428572
$ws.Range("M16").Value = 0.02285714285714286
$ws.Range("N16").Value = 0.005714285714285714
$ws.Range("O16").Value = 0.05714285714285714
$ws.Range("S16").Value = 0.12
$ws.Range("F17").Value = 0.02102803738317757
$ws.Range("H17").Value = 0.1845794392523364
$ws.Range("I17").Value = 0.1004672897196262
$ws.Range("J17").Value = 0.4275700934579439
$ws.Range("K17").Value = 0.08177570093457943
$ws.Range("M17").Value = 0.01401869158878505
$ws.Range("N17").Value = 0.002336448598130841
$ws.Range("O17").Value = 0.06775700934579439
$ws.Range("S17").Value = 0.1004672897196262
$ws.Range("F18").Value = 0.0125
$ws.Range("H18").Value = 0.15
$ws.Range("I18").Value = 0.1125
$ws.Range("J18").Value = 0.3875
$ws.Range("K18").Value = 0.1166666666666667
$ws.Range("M18").Value = 0.025
$ws.Range("O18").Value = 0.04583333333333333
$ws.Range("S18").Value = 0.15
$ws.Range("F19").Value = 0.01240875912408759
$ws.Range("H19").Value = 0.2102189781021898
$ws.Range("I19").Value = 0.0948905109489051
$ws.Range("J19").Value = 0.3408759124087591
$ws.Range("K19").Value = 0.1021897810218978
$ws.Range("M19").Value = 0.01970802919708029
$ws.Range("O19").Value = 0.06715328467153285
$ws.Range("S19").Value = 0.1525547445255475
